# Thpo-Mpl LR-pair sheet refresh with new TPM-derived NATMI output.
#
# The re-run of the pipeline with the new TPM numbers adds a new
# "ECs -> Thpo -> Mpl -> ECs" sending-cluster row and recalculates the
# receptor-side specificity figures for every row (the receptor is now
# seen as expressed by 2 of the 3 clusters instead of just 1, which
# changes every "receptor ..." / "edge ..." derived column).
#
# Resulting layout (row 1 is the existing header row, untouched):
#   row 2: ECs   -> Thpo -> Mpl -> ECs
#   row 3: FAPs  -> Thpo -> Mpl -> ECs
#   row 4: MuSCs -> Thpo -> Mpl -> ECs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: ECs -> Thpo -> Mpl -> ECs -------------------------------
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Thpo"
$ws.Range("C2").Value = "Mpl"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.001845666666666667
$ws.Range("H2").Value = 0.005537
$ws.Range("I2").Value = 0.0009892926977084571
$ws.Range("J2").Value = 0.0009892926977084571
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1272896666666667
$ws.Range("N2").Value = 0.381869
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0002349342947777778
$ws.Range("R2").Value = 0.002114408653
$ws.Range("S2").Value = 0.0009892926977084571
$ws.Range("T2").Value = 0.0009892926977084571

# ---- Row 3: FAPs -> Thpo -> Mpl -> ECs (was row 2, recalculated) ----
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Thpo"
$ws.Range("C3").Value = "Mpl"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.737773666666667
$ws.Range("H3").Value = 5.213321000000001
$ws.Range("I3").Value = 0.931461151545991
$ws.Range("J3").Value = 0.9314611515459909
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1272896666666667
$ws.Range("N3").Value = 0.381869
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.2212006307721111
$ws.Range("R3").Value = 1.990805676949
$ws.Range("S3").Value = 0.931461151545991
$ws.Range("T3").Value = 0.9314611515459909

# ---- Row 4: MuSCs -> Thpo -> Mpl -> ECs (was row 3, recalculated) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Thpo"
$ws.Range("C4").Value = "Mpl"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1260233333333333
$ws.Range("H4").Value = 0.37807
$ws.Range("I4").Value = 0.0675495557563006
$ws.Range("J4").Value = 0.06754955575630059
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1272896666666667
$ws.Range("N4").Value = 0.381869
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.01604146809222222
$ws.Range("R4").Value = 0.14437321283
$ws.Range("S4").Value = 0.0675495557563006
$ws.Range("T4").Value = 0.06754955575630059
